$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.508.19"
$ws.Range("E2").Value = "  +5.64%  "
$ws.Range("D3").Value = "1.722.39"
$ws.Range("E3").Value = "  +4.42%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "226.08"
$ws.Range("E5").Value = "  +3.60%  "
$ws.Range("D6").Value = "0.5338"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "0.2661"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("D9").Value = "0.06584"
$ws.Range("E9").Value = "  +4.22%  "
$ws.Range("D10").Value = "21.66"
$ws.Range("E10").Value = "  +6.36%  "
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").Value = "4.619"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "1.725.40"
$ws.Range("E13").Value = "  +4.63%  "
$ws.Range("D14").Value = "1.960.00"
$ws.Range("E14").Value = "  +4.40%  "
$ws.Range("D15").Value = "0.5827"
$ws.Range("E15").Value = "  +4.29%  "
$ws.Range("D16").Value = "0.0₅8291"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "67.79"
$ws.Range("E17").Value = "  +4.00%  "
$ws.Range("D18").Value = "27.516.88"
$ws.Range("E18").Value = "  +5.74%  "
$ws.Range("D19").Value = "219.62"
$ws.Range("E19").Value = "  +14.71%  "
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "4.728"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").Value = "10.63"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("D23").Value = "6.073"
$ws.Range("E23").Value = "  +2.95%  "
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "147.89"
$ws.Range("E25").Value = "  +3.07%  "
$ws.Range("D26").Value = "1.740"
$ws.Range("E26").Value = "  +15.42%  "
$ws.Range("D27").Value = "0.1233"
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("D28").Value = "7.402"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").Value = "16.54"
$ws.Range("E29").Value = "  +4.33%  "
$ws.Range("D30").Value = "0.05550"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("E31").Value = "  +2.73%  "
$ws.Range("D32").Value = "3.561"
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D33").Value = "3.442"
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("D34").Value = "1.661"
$ws.Range("E34").Value = "  +6.99%  "
$ws.Range("D35").Value = "2.862"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").Value = "0.9634"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").Value = "2.424"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "0.5956"
$ws.Range("E38").Value = "  +5.79%  "
$ws.Range("E39").Value = "  +4.69%  "
$ws.Range("D40").Value = "5.909"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").Value = "0.8546"
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("D42").Value = "1.055.96"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").Value = "1.006"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "101.37"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "1.866.62"
$ws.Range("E45").Value = "  +4.46%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  +4.26%  "
$ws.Range("D47").Value = "58.89"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").Value = "8.199"
$ws.Range("E48").Value = "  +3.63%  "
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").Value = "0.05245"
$ws.Range("E51").Value = "  +1.97%  "
